$wb = $excel.ActiveWorkbook

# --- Sheet: ASCOM ItelescopeV3 Methods ---
$wsMethods = $wb.Worksheets.Item("ASCOM ItelescopeV3 Methods")

# PulseGuide method row (row 14) is now implemented:
#   Supported column -> "Yes"
#   Return value column -> "Nothing"
# Copy from existing cells holding the same literal text so the values are
# written as plain text (and not re-interpreted, e.g. "True"/"False" words
# as booleans) exactly like the rest of the sheet.
[void]$wsMethods.Range("B2").Copy($wsMethods.Range("B14"))
[void]$wsMethods.Range("D6").Copy($wsMethods.Range("D14"))

# --- Sheet: ASCOM ItelescopeV3 Properties ---
$wsProps = $wb.Worksheets.Item("ASCOM ItelescopeV3 Properties")

# CanPulseGuide property row (row 11) is now implemented -> "True"
[void]$wsProps.Range("D10").Copy($wsProps.Range("D11"))

# IsPulseGuiding property row (row 37) is now implemented:
#   Supported column -> "Yes"
#   Return type column -> "True or False"
[void]$wsProps.Range("B2").Copy($wsProps.Range("B37"))
[void]$wsProps.Range("D7").Copy($wsProps.Range("D37"))

# Update the visible selection/scroll position on the Properties sheet (it is not the
# active tab, so select it, adjust the view, then return focus to the Methods sheet).
$wsProps.Activate()
[void]$wsProps.Range("A37").Select()
$excel.ActiveWindow.ScrollRow = 22

# Update the visible selection on the Methods sheet and leave it as the active tab
$wsMethods.Activate()
[void]$wsMethods.Range("B15").Select()

# Restore the application window position/size recorded in the workbook view
# (best-effort: some hosts do not persist OS window chrome back to the file).
$excel.Left = 28680
$excel.Top = -75
